$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; existing rows 4..30 shift down to 5..31.
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the new weekly record.
$ws.Cells.Item(4, 1).Value = 9
$ws.Cells.Item(4, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = Get-Date -Year 2022 -Month 5 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = 100112035
$ws.Cells.Item(4, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(4, 11).Value = 25000
$ws.Cells.Item(4, 12).Value = 25000
$ws.Cells.Item(4, 13).Value = 25000
$ws.Cells.Item(4, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(4, 15).Value = "Hijuelas"
$ws.Cells.Item(4, 16).Value = 1667
$ws.Cells.Item(4, 17).Value = 15
$ws.Cells.Item(4, 18).Value = "Hortaliza"

Write-Output "done"
